# "custom functions using chatgpt in lab apps removed"
#
# This removes the custom "set_impression_of_dialogue(&impression)" action calls
# (column G, the "actions" column) from the #final_with_impression related
# transitions, and replaces the hand-written "{impression} . Thank you for
# your time." system utterance with an inline ChatGPT-prompt utterance
# "{$"Generate a short utterance to say the system's impression."} Thank you
# for your time." so the impression is generated directly in the utterance
# instead of via a separate custom action/variable.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Remove the "set_impression_of_dialogue(&impression)" action cells entirely
# (not just clear their contents) for rows 20, 21 and 22 (new-sandwich,
# egg-salad-sandwich and known-sandwich states).
$ws.Range("G20").Clear()
$ws.Range("G21").Clear()
$ws.Range("G22").Clear()

# Replace the #final_with_impression system utterance so it generates the
# impression inline via a ChatGPT prompt instead of relying on the removed
# custom function/variable.
$ws.Range("C25").Value = "{$""Generate a short utterance to say the system's impression.""} Thank you for your time."

# Match the author's final cursor position after the edit.
$ws.Range("F20").Select()
